$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 202; this shifts existing rows 202-240 down to 203-241
# and Excel auto-extends the used-range dimension (A1:R240 -> A1:R241).
$ws.Rows.Item(202).EntireRow.Insert()

# Populate the newly inserted row 202 with its data.
$ws.Range("A202").Value = 4
$ws.Range("B202").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C202").Value = "Los Lagos"
$ws.Range("D202").Value = 44711
$ws.Range("E202").Value = 10
$ws.Range("F202").Value = 100112032
$ws.Range("G202").Value = "Zapallo italiano"
$ws.Range("H202").Value = "Sin especificar"
$ws.Range("I202").Value = "Primera"
$ws.Range("J202").Value = 70
$ws.Range("K202").Value = 21000
$ws.Range("L202").Value = 21000
$ws.Range("M202").Value = 21000
$ws.Range("N202").Value = "`$/caja 50 unidades"
$ws.Range("O202").Value = "Región de Arica y Parinacota"
$ws.Range("P202").Value = 420
$ws.Range("Q202").Value = 50
$ws.Range("R202").Value = "Hortaliza"
